$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Il6"
$ws.Cells.Item(2,3).Value = "Il6ra"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 7.292895000000001
$ws.Cells.Item(2,8).Value = 21.878685
$ws.Cells.Item(2,9).Value = 0.2565758520803378
$ws.Cells.Item(2,10).Value = 0.2565758520803378
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 13.338276
$ws.Cells.Item(2,14).Value = 40.01482799999999
$ws.Cells.Item(2,15).Value = 0.80809692568033
$ws.Cells.Item(2,16).Value = 0.8080969256803301
$ws.Cells.Item(2,17).Value = 97.27464634902
$ws.Cells.Item(2,18).Value = 875.4718171411799
$ws.Cells.Item(2,19).Value = 0.2073381572699321
$ws.Cells.Item(2,20).Value = 0.2073381572699321

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Il6"
$ws.Cells.Item(3,3).Value = "Il6ra"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 7.292895000000001
$ws.Cells.Item(3,8).Value = 21.878685
$ws.Cells.Item(3,9).Value = 0.2565758520803378
$ws.Cells.Item(3,10).Value = 0.2565758520803378
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.961838
$ws.Cells.Item(3,14).Value = 8.885514000000001
$ws.Cells.Item(3,15).Value = 0.1794423943666466
$ws.Cells.Item(3,16).Value = 0.1794423943666466
$ws.Cells.Item(3,17).Value = 21.60037354101
$ws.Cells.Item(3,18).Value = 194.40336186909
$ws.Cells.Item(3,19).Value = 0.04604058523395836
$ws.Cells.Item(3,20).Value = 0.04604058523395837

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Il6"
$ws.Cells.Item(4,3).Value = "Il6ra"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 7.292895000000001
$ws.Cells.Item(4,8).Value = 21.878685
$ws.Cells.Item(4,9).Value = 0.2565758520803378
$ws.Cells.Item(4,10).Value = 0.2565758520803378
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.2056733333333333
$ws.Cells.Item(4,14).Value = 0.6170199999999999
$ws.Cells.Item(4,15).Value = 0.01246067995302334
$ws.Cells.Item(4,16).Value = 0.01246067995302335
$ws.Cells.Item(4,17).Value = 1.4999540243
$ws.Cells.Item(4,18).Value = 13.4995862187
$ws.Cells.Item(4,19).Value = 0.003197109576447348
$ws.Cells.Item(4,20).Value = 0.003197109576447349

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Il6"
$ws.Cells.Item(5,3).Value = "Il6ra"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 14.46476266666667
$ws.Cells.Item(5,8).Value = 43.394288
$ws.Cells.Item(5,9).Value = 0.5088937666509471
$ws.Cells.Item(5,10).Value = 0.5088937666509471
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 13.338276
$ws.Cells.Item(5,14).Value = 40.01482799999999
$ws.Cells.Item(5,15).Value = 0.80809692568033
$ws.Cells.Item(5,16).Value = 0.8080969256803301
$ws.Cells.Item(5,17).Value = 192.934996722496
$ws.Cells.Item(5,18).Value = 1736.414970502464
$ws.Cells.Item(5,19).Value = 0.4112354883285136
$ws.Cells.Item(5,20).Value = 0.4112354883285136

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Il6"
$ws.Cells.Item(6,3).Value = "Il6ra"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 14.46476266666667
$ws.Cells.Item(6,8).Value = 43.394288
$ws.Cells.Item(6,9).Value = 0.5088937666509471
$ws.Cells.Item(6,10).Value = 0.5088937666509471
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.961838
$ws.Cells.Item(6,14).Value = 8.885514000000001
$ws.Cells.Item(6,15).Value = 0.1794423943666466
$ws.Cells.Item(6,16).Value = 0.1794423943666466
$ws.Cells.Item(6,17).Value = 42.84228372711467
$ws.Cells.Item(6,18).Value = 385.5805535440321
$ws.Cells.Item(6,19).Value = 0.09131711596610748
$ws.Cells.Item(6,20).Value = 0.09131711596610749

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Il6"
$ws.Cells.Item(7,3).Value = "Il6ra"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 14.46476266666667
$ws.Cells.Item(7,8).Value = 43.394288
$ws.Cells.Item(7,9).Value = 0.5088937666509471
$ws.Cells.Item(7,10).Value = 0.5088937666509471
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.2056733333333333
$ws.Cells.Item(7,14).Value = 0.6170199999999999
$ws.Cells.Item(7,15).Value = 0.01246067995302334
$ws.Cells.Item(7,16).Value = 0.01246067995302335
$ws.Cells.Item(7,17).Value = 2.975015953528888
$ws.Cells.Item(7,18).Value = 26.77514358176
$ws.Cells.Item(7,19).Value = 0.006341162356325996
$ws.Cells.Item(7,20).Value = 0.006341162356325997

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Il6"
$ws.Cells.Item(8,3).Value = "Il6ra"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 6.666276
$ws.Cells.Item(8,8).Value = 19.998828
$ws.Cells.Item(8,9).Value = 0.2345303812687151
$ws.Cells.Item(8,10).Value = 0.2345303812687151
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 13.338276
$ws.Cells.Item(8,14).Value = 40.01482799999999
$ws.Cells.Item(8,15).Value = 0.80809692568033
$ws.Cells.Item(8,16).Value = 0.8080969256803301
$ws.Cells.Item(8,17).Value = 88.91662918017599
$ws.Cells.Item(8,18).Value = 800.2496626215839
$ws.Cells.Item(8,19).Value = 0.1895232800818843
$ws.Cells.Item(8,20).Value = 0.1895232800818843

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Il6"
$ws.Cells.Item(9,3).Value = "Il6ra"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 6.666276
$ws.Cells.Item(9,8).Value = 19.998828
$ws.Cells.Item(9,9).Value = 0.2345303812687151
$ws.Cells.Item(9,10).Value = 0.2345303812687151
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.961838
$ws.Cells.Item(9,14).Value = 8.885514000000001
$ws.Cells.Item(9,15).Value = 0.1794423943666466
$ws.Cells.Item(9,16).Value = 0.1794423943666466
$ws.Cells.Item(9,17).Value = 19.744429575288
$ws.Cells.Item(9,18).Value = 177.699866177592
$ws.Cells.Item(9,19).Value = 0.04208469316658076
$ws.Cells.Item(9,20).Value = 0.04208469316658076

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Il6"
$ws.Cells.Item(10,3).Value = "Il6ra"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 6.666276
$ws.Cells.Item(10,8).Value = 19.998828
$ws.Cells.Item(10,9).Value = 0.2345303812687151
$ws.Cells.Item(10,10).Value = 0.2345303812687151
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.2056733333333333
$ws.Cells.Item(10,14).Value = 0.6170199999999999
$ws.Cells.Item(10,15).Value = 0.01246067995302334
$ws.Cells.Item(10,16).Value = 0.01246067995302335
$ws.Cells.Item(10,17).Value = 1.37107520584
$ws.Cells.Item(10,18).Value = 12.33967685256
$ws.Cells.Item(10,19).Value = 0.00292240802025
$ws.Cells.Item(10,20).Value = 0.00292240802025
